# Reorders / relabels the appointment schedule table in RUNA COFFEE.docx.
# The table keeps the same number of rows; only the FRANJA HORARIA (col 1)
# and COMPRADOR (col 3) text contents are reassigned per row. Column 2
# (MESA) remains empty in every row and is left untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row index => (new time slot, new buyer)
$updates = @(
    @{ Row = 2; Time = "08:30 - 08:45"; Buyer = "FLOR A FRUTO" },
    @{ Row = 3; Time = "08:45 - 09:00"; Buyer = "INMERSSO BOUTIQUE" },
    @{ Row = 4; Time = "09:15 - 09:30"; Buyer = "BOX BRAND" },
    @{ Row = 5; Time = "09:45 - 10:00"; Buyer = "CAFÉ MOLINA" },
    @{ Row = 6; Time = "10:00 - 10:15"; Buyer = "COLFRESH COFFEE" },
    @{ Row = 7; Time = "10:15 - 10:30"; Buyer = "NEIRA YORK COFFEE" },
    @{ Row = 9; Time = "10:45 - 11:00"; Buyer = "ARMANDO VELÁSQUEZ" }
)

foreach ($u in $updates) {
    $t.Cell($u.Row, 1).Range.Text = $u.Time
    $t.Cell($u.Row, 3).Range.Text = $u.Buyer
}
